$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("测试概览")
$ws1.Columns.Item(2).Delete()

$ws2 = $wb.Worksheets.Item("BUG汇总 ")
$ws2.Columns.Item(4).Delete()
